$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 301.2
$ws.Range("I5").Value = 88.42856999999999
$ws.Range("J5").Value = 487.375
$ws.Range("K5").Value = 88.42856999999999
$ws.Range("L5").Value = 487.375
$ws.Range("M5").Value = 26.57143000000001
$ws.Range("N5").Value = -717.375

$ws.Range("H62").Value = 27781262
$ws.Range("I62").Value = 41670484
$ws.Range("J62").Value = 2820
$ws.Range("K62").Value = 41670484
$ws.Range("L62").Value = 2820
$ws.Range("M62").Value = -41669860
$ws.Range("N62").Value = -4068

$ws.Range("H65").Value = 27781262
$ws.Range("I65").Value = 41670484
$ws.Range("J65").Value = 2820
$ws.Range("K65").Value = 208352420
$ws.Range("L65").Value = 14100
$ws.Range("M65").Value = -208349300
$ws.Range("N65").Value = -20340

$ws.Range("H88").Value = 9273.267
$ws.Range("I88").Value = 2482.8333
$ws.Range("J88").Value = 13800.223
$ws.Range("K88").Value = 2482.8333
$ws.Range("L88").Value = 13800.223
$ws.Range("M88").Value = -2076.8333
$ws.Range("N88").Value = -14612.223

$ws.Range("H91").Value = 9273.267
$ws.Range("I91").Value = 2482.8333
$ws.Range("J91").Value = 13800.223
$ws.Range("K91").Value = 2482.8333
$ws.Range("L91").Value = 13800.223
$ws.Range("M91").Value = -1078.8333
$ws.Range("N91").Value = -16608.223

$ws.Range("H132").Value = 1964.8154
$ws.Range("I132").Value = 1344.3818
$ws.Range("J132").Value = 5377.2
$ws.Range("K132").Value = 4033.1454
$ws.Range("L132").Value = 16131.6
$ws.Range("M132").Value = -1503.1454
$ws.Range("N132").Value = -21191.6

$ws.Range("H135").Value = 2271.0312
$ws.Range("I135").Value = 2077.2856
$ws.Range("J135").Value = 2640.9092
$ws.Range("K135").Value = 18695.5704
$ws.Range("L135").Value = 23768.1828
$ws.Range("M135").Value = -16160.5704
$ws.Range("N135").Value = -28838.1828

$ws.Range("H138").Value = 2151.0603
$ws.Range("I138").Value = 1234
$ws.Range("J138").Value = 2888.6956
$ws.Range("K138").Value = 3702
$ws.Range("L138").Value = 8666.086800000001
$ws.Range("M138").Value = 1438
$ws.Range("N138").Value = -18946.0868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1055
$ws.Range("I2").Value = 1064
$ws.Range("J2").Value = 1028
$ws.Range("K2").Value = 1064
$ws.Range("L2").Value = 1028
$ws.Range("M2").Value = -951
$ws.Range("N2").Value = -1254

$ws.Range("H74").Value = 1849.1333
$ws.Range("I74").Value = 1402.3
$ws.Range("J74").Value = 2742.8
$ws.Range("K74").Value = 1402.3
$ws.Range("L74").Value = 2742.8
$ws.Range("M74").Value = -528.3
$ws.Range("N74").Value = -4490.8

$ws.Range("H77").Value = 1849.1333
$ws.Range("I77").Value = 1402.3
$ws.Range("J77").Value = 2742.8
$ws.Range("K77").Value = 7011.5
$ws.Range("L77").Value = 13714
$ws.Range("M77").Value = -2643.5
$ws.Range("N77").Value = -22450

$ws.Range("H116").Value = 1055
$ws.Range("I116").Value = 1064
$ws.Range("J116").Value = 1028
$ws.Range("K116").Value = 1064
$ws.Range("L116").Value = 1028
$ws.Range("M116").Value = 1230
$ws.Range("N116").Value = -5616

$ws.Range("H122").Value = 612669.4
$ws.Range("I122").Value = 803357.1
$ws.Range("J122").Value = 2468.6
$ws.Range("K122").Value = 2410071.3
$ws.Range("L122").Value = 7405.799999999999
$ws.Range("M122").Value = -2407621.3
$ws.Range("N122").Value = -12305.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1055
$ws.Range("I3").Value = 1064
$ws.Range("J3").Value = 1028
$ws.Range("K3").Value = 1064
$ws.Range("L3").Value = 1028
$ws.Range("M3").Value = -950
$ws.Range("N3").Value = -1256

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 49613.863
$ws.Range("I19").Value = 250376.25
$ws.Range("K19").Value = 250376.25
$ws.Range("M19").Value = -250206.25

$ws.Range("H24").Value = 49613.863
$ws.Range("I24").Value = 250376.25
$ws.Range("K24").Value = 250376.25
$ws.Range("M24").Value = -250206.25

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 8369.429
$ws.Range("J94").Value = 8431.833000000001
$ws.Range("L94").Value = 8431.833000000001
$ws.Range("N94").Value = -9333.833000000001

$ws.Range("H132").Value = 3071.75
$ws.Range("I132").Value = 2437
$ws.Range("J132").Value = 3960.4
$ws.Range("K132").Value = 7311
$ws.Range("L132").Value = 11881.2
$ws.Range("M132").Value = -4781
$ws.Range("N132").Value = -16941.2

$ws.Range("H141").Value = 272494
$ws.Range("I141").Value = 39811
$ws.Range("J141").Value = 291884.25
$ws.Range("K141").Value = 39811
$ws.Range("L141").Value = 291884.25
$ws.Range("M141").Value = -34631
$ws.Range("N141").Value = -302244.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5122.815
$ws.Range("I5").Value = 7151.067
$ws.Range("J5").Value = 2587.5
$ws.Range("K5").Value = 21453.201
$ws.Range("L5").Value = 7762.5
$ws.Range("M5").Value = -21341.201
$ws.Range("N5").Value = -7986.5

$ws.Range("H129").Value = 1385.4
$ws.Range("I129").Value = 791.6667
$ws.Range("J129").Value = 1781.2222
$ws.Range("K129").Value = 2375.0001
$ws.Range("L129").Value = 5343.6666
$ws.Range("M129").Value = 2624.9999
$ws.Range("N129").Value = -15343.6666

$ws.Range("H135").Value = 5122.815
$ws.Range("I135").Value = 7151.067
$ws.Range("J135").Value = 2587.5
$ws.Range("K135").Value = 64359.603
$ws.Range("L135").Value = 23287.5
$ws.Range("M135").Value = -61824.603
$ws.Range("N135").Value = -28357.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 59160748
$ws.Range("I122").Value = 96803950
$ws.Range("J122").Value = 7142.2856
$ws.Range("K122").Value = 290411850
$ws.Range("L122").Value = 21426.8568
$ws.Range("M122").Value = -290409400
$ws.Range("N122").Value = -26326.8568

$ws.Range("H123").Value = 28601
$ws.Range("J123").Value = 28601
$ws.Range("L123").Value = 28601
$ws.Range("N123").Value = -33501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 2799.375
$ws.Range("J38").Value = 2799.375
$ws.Range("L38").Value = 2799.375
$ws.Range("N38").Value = -3619.375

$ws.Range("H122").Value = 5095960
$ws.Range("I122").Value = 5960413.5
$ws.Range("J122").Value = 2502600
$ws.Range("K122").Value = 17881240.5
$ws.Range("L122").Value = 7507800
$ws.Range("M122").Value = -17878790.5
$ws.Range("N122").Value = -7512700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 7210
$ws.Range("J30").Value = 7210
$ws.Range("L30").Value = 7210
$ws.Range("N30").Value = -7424

$ws.Range("H104").Value = 36891.43
$ws.Range("J104").Value = 36891.43
$ws.Range("L104").Value = 36891.43
$ws.Range("N104").Value = -43879.43
